$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.905.38"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.541.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.23"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.25"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0819"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.63"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.110"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.929.97"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.545.03"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.852"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.970.96"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.88"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.79"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.77"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.94"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.50"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.39%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  +3.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "41.01"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.28%  "
$ws.Range("E30").Value = "  +4.23%  "
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.33"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.17"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.38"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("E35").Value = "  +4.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.96"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0794"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("E39").Value = "  +15.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.118"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.94"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -10.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.84"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.023.72"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.40"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.88"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "106.58"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.784.02"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.91%  "
